$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 480, shifting existing rows 480:505 down to 482:507.
$ws.Rows.Item(480).Resize(2).Insert()

# Fill the two newly inserted rows (480 and 481) with the new data rows.
$ws.Range("A480").Value = 6
$ws.Range("B480").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C480").Value = "Metropolitana"
$ws.Range("D480").Value = 44931
$ws.Range("E480").Value = 13
$ws.Range("F480").Value = "Fruta"
$ws.Range("G480").Value = 100101
$ws.Range("H480").Value = "Berries"
$ws.Range("I480").Value = 100101001
$ws.Range("J480").Value = "Arándano (blue)"
$ws.Range("K480").Value = "Sin especificar"
$ws.Range("L480").Value = "Especial"
$ws.Range("M480").Value = 250
$ws.Range("N480").Value = 3000
$ws.Range("O480").Value = 3000
$ws.Range("P480").Value = 3000
$ws.Range("Q480").Value = "$/bandeja 2 kilos"
$ws.Range("R480").Value = "Provincia de Curicó"
$ws.Range("S480").Value = 1500
$ws.Range("T480").Value = 2

$ws.Range("A481").Value = 6
$ws.Range("B481").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C481").Value = "Metropolitana"
$ws.Range("D481").Value = 44931
$ws.Range("E481").Value = 13
$ws.Range("F481").Value = "Fruta"
$ws.Range("G481").Value = 100101
$ws.Range("H481").Value = "Berries"
$ws.Range("I481").Value = 100101001
$ws.Range("J481").Value = "Arándano (blue)"
$ws.Range("K481").Value = "Sin especificar"
$ws.Range("L481").Value = "Especial"
$ws.Range("M481").Value = 2000
$ws.Range("N481").Value = 3000
$ws.Range("O481").Value = 3000
$ws.Range("P481").Value = 3000
$ws.Range("Q481").Value = "$/bandeja 2 kilos"
$ws.Range("R481").Value = "Región del Maule"
$ws.Range("S481").Value = 1500
$ws.Range("T481").Value = 2

# Match the date cell style/format used by the rest of column D.
$ws.Range("D480:D481").NumberFormat = "YYYY-MM-DD HH:MM:SS"
